$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.059.13"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "2.360.40"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.972"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.57%  "
$ws.Range("D15").Value = "2.721.19"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "2.357.36"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "45.068.20"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.73%  "
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.30%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0974"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "167.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +4.61%  "
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.71%  "
$ws.Range("D47").Value = "1.840.40"
$ws.Range("E47").Value = "  +10.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.17%  "
$ws.Range("E49").Value = "  +6.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.53%  "
